# Regenerate save_data: replace column G ("K", previously derived from
# Strike#) with freshly calculated s_vals for rows 2-34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 5
    4  = 0
    5  = 2
    6  = 0
    7  = 1
    8  = 5
    9  = 2
    10 = 2
    11 = 1
    12 = 3
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 2
    18 = 2
    19 = 1
    20 = 0
    21 = 1
    22 = 0
    23 = 1
    24 = 0
    25 = 0
    26 = 2
    27 = 0
    28 = 2
    29 = 1
    30 = 3
    31 = 6
    32 = 1
    33 = 1
    34 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
